# Edit script: updates Rhode Island converted data sheet
# - Row 7: T7/U7 set to 0, AA7 set to 12
# - Rows 21-221: AA column (lockdown effectiveness score) recalculated for a 24-day window
# - Rows 222-233: new date rows appended (9/30/2020 - 10/11/2020), copied policy indicator
#   pattern from row 221, with AA = 0.4166666666666667

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7 updates ---
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 0
$ws.Cells.Item(7, 27).Value = 12

# --- AA column (lockdown effectiveness score) updates for rows 21-221 ---
$aaUpdates = @(
    "21:0.04166666666666666",
    "22:0.04166666666666666",
    "23:0.04166666666666666",
    "24:0.125",
    "25:0.2916666666666667",
    "26:0.2916666666666667",
    "27:0.2916666666666667",
    "28:0.2916666666666667",
    "29:0.2916666666666667",
    "30:0.2916666666666667",
    "31:0.4166666666666667",
    "32:0.4166666666666667",
    "33:0.4166666666666667",
    "34:0.4166666666666667",
    "35:0.4166666666666667",
    "36:0.6666666666666666",
    "37:0.6666666666666666",
    "38:0.75",
    "39:0.75",
    "40:0.75",
    "41:0.75",
    "42:0.8333333333250001",
    "43:0.8333333333250001",
    "44:0.8333333333250001",
    "45:0.8333333333250001",
    "46:0.8333333333250001",
    "47:0.8333333333250001",
    "48:0.8333333333250001",
    "49:0.8333333333250001",
    "50:0.8333333333250001",
    "51:0.8333333333250001",
    "52:0.8333333333250001",
    "53:0.8333333333250001",
    "54:0.8333333333250001",
    "55:0.8333333333250001",
    "56:0.8333333333250001",
    "57:0.8333333333250001",
    "58:0.8333333333250001",
    "59:0.8333333333250001",
    "60:0.8333333333250001",
    "61:0.8333333333250001",
    "62:0.8333333333250001",
    "63:0.8333333333250001",
    "64:0.8333333333250001",
    "65:0.8333333333250001",
    "66:0.8333333333250001",
    "67:0.8333333333250001",
    "68:0.8333333333250001",
    "69:0.8333333333250001",
    "70:0.8333333333250001",
    "71:0.8333333333250001",
    "72:0.8333333333250001",
    "73:0.8333333333250001",
    "74:0.8333333333250001",
    "75:0.8333333333250001",
    "76:0.8333333333250001",
    "77:0.9166666666583335",
    "78:0.8333333333250001",
    "79:0.6888888888833334",
    "80:0.6888888888833334",
    "81:0.6888888888833334",
    "82:0.6888888888833334",
    "83:0.6888888888833334",
    "84:0.6888888888833334",
    "85:0.6888888888833334",
    "86:0.6888888888833334",
    "87:0.6888888888833334",
    "88:0.6888888888833334",
    "89:0.6888888888833334",
    "90:0.6888888888833334",
    "91:0.6888888888833334",
    "92:0.6888888888833334",
    "93:0.6888888888833334",
    "94:0.6888888888833334",
    "95:0.6888888888833334",
    "96:0.6888888888833334",
    "97:0.6888888888833334",
    "98:0.6888888888833334",
    "99:0.6888888888833334",
    "100:0.6888888888833334",
    "101:0.6888888888833334",
    "102:0.4861111111083334",
    "103:0.4861111111083334",
    "104:0.4861111111083334",
    "105:0.4861111111083334",
    "106:0.4861111111083334",
    "107:0.4861111111083334",
    "108:0.4861111111083334",
    "109:0.4444444444416667",
    "110:0.4444444444416667",
    "111:0.4444444444416667",
    "112:0.4444444444416667",
    "113:0.4444444444416667",
    "114:0.4444444444416667",
    "115:0.4444444444416667",
    "116:0.4444444444416667",
    "117:0.4444444444416667",
    "118:0.4444444444416667",
    "119:0.4444444444416667",
    "120:0.4444444444416667",
    "121:0.4444444444416667",
    "122:0.4444444444416667",
    "123:0.4444444444416667",
    "124:0.4444444444416667",
    "125:0.4444444444416667",
    "126:0.4444444444416667",
    "127:0.4444444444416667",
    "128:0.4444444444416667",
    "129:0.4444444444416667",
    "130:0.4444444444416667",
    "131:0.4166666666666667",
    "132:0.4166666666666667",
    "133:0.4166666666666667",
    "134:0.4166666666666667",
    "135:0.4166666666666667",
    "136:0.4166666666666667",
    "137:0.4166666666666667",
    "138:0.4166666666666667",
    "139:0.4166666666666667",
    "140:0.4166666666666667",
    "141:0.4166666666666667",
    "142:0.4166666666666667",
    "143:0.4166666666666667",
    "144:0.4166666666666667",
    "145:0.4166666666666667",
    "146:0.4166666666666667",
    "147:0.4166666666666667",
    "148:0.4166666666666667",
    "149:0.4166666666666667",
    "150:0.4166666666666667",
    "151:0.4166666666666667",
    "152:0.4166666666666667",
    "153:0.4166666666666667",
    "154:0.4166666666666667",
    "155:0.4166666666666667",
    "156:0.4166666666666667",
    "157:0.4166666666666667",
    "158:0.4166666666666667",
    "159:0.4166666666666667",
    "160:0.4166666666666667",
    "161:0.4166666666666667",
    "162:0.4166666666666667",
    "163:0.4166666666666667",
    "164:0.4166666666666667",
    "165:0.4166666666666667",
    "166:0.4166666666666667",
    "167:0.4166666666666667",
    "168:0.4166666666666667",
    "169:0.4166666666666667",
    "170:0.4166666666666667",
    "171:0.4166666666666667",
    "172:0.4166666666666667",
    "173:0.4166666666666667",
    "174:0.4166666666666667",
    "175:0.4166666666666667",
    "176:0.4166666666666667",
    "177:0.4166666666666667",
    "178:0.4166666666666667",
    "179:0.4166666666666667",
    "180:0.4166666666666667",
    "181:0.4166666666666667",
    "182:0.4166666666666667",
    "183:0.4166666666666667",
    "184:0.4166666666666667",
    "185:0.4166666666666667",
    "186:0.4166666666666667",
    "187:0.4166666666666667",
    "188:0.4166666666666667",
    "189:0.4166666666666667",
    "190:0.4166666666666667",
    "191:0.4166666666666667",
    "192:0.4166666666666667",
    "193:0.4166666666666667",
    "194:0.4166666666666667",
    "195:0.4166666666666667",
    "196:0.4166666666666667",
    "197:0.4166666666666667",
    "198:0.4166666666666667",
    "199:0.4166666666666667",
    "200:0.4166666666666667",
    "201:0.4166666666666667",
    "202:0.4166666666666667",
    "203:0.4166666666666667",
    "204:0.4166666666666667",
    "205:0.4166666666666667",
    "206:0.4166666666666667",
    "207:0.4166666666666667",
    "208:0.4166666666666667",
    "209:0.4166666666666667",
    "210:0.4166666666666667",
    "211:0.4166666666666667",
    "212:0.4166666666666667",
    "213:0.4166666666666667",
    "214:0.4166666666666667",
    "215:0.4166666666666667",
    "216:0.4166666666666667",
    "217:0.4166666666666667",
    "218:0.4166666666666667",
    "219:0.4166666666666667",
    "220:0.4166666666666667",
    "221:0.4166666666666667"
)
foreach ($item in $aaUpdates) {
    $parts = $item.Split(":")
    $row = [int]$parts[0]
    $val = [double]$parts[1]
    $ws.Cells.Item($row, 27).Value = $val
}

# --- New rows 222-233: additional dates 9/30/2020 - 10/11/2020 ---
$newDates = @(
    "9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020",
    "10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020"
)

# Values for columns B..Z are identical to row 221 (copied forward), AA is the new score
$rowValues = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1)
$aaValue = 0.4166666666666667

$srcLabelCell = $ws.Cells.Item(221, 1)

$newRow = 222
foreach ($d in $newDates) {
    $labelCell = $ws.Cells.Item($newRow, 1)
    # Set value as text first (apostrophe prefix forces text, avoiding Excel's
    # automatic date parsing), then paste the source cell's formatting on top
    # so the label keeps the same bold/border/text style used by every other
    # date label in column A (style index 1).
    $labelCell.Value = "'" + $d
    $srcLabelCell.Copy()
    $labelCell.PasteSpecial(-4122)

    for ($col = 2; $col -le 26; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $rowValues[$col - 2]
    }
    $ws.Cells.Item($newRow, 27).Value = $aaValue

    $newRow = $newRow + 1
}

$excel.CutCopyMode = 0

Write-Host "Edit complete"
